# Auto-generated edit script applying the country/provincias Spain daily data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 00:13"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6630867
$ws.Range("C4").Value = 41220
$ws.Range("D4").Value = 3894339
$ws.Range("E4").Value = 2539280
$ws.Range("G4").Value = 921
$ws.Range("H4").Value = 197248

# Row 6: Brasil
$ws.Range("B6").Value = 4282164
$ws.Range("C6").Value = 42401
$ws.Range("E6").Value = 654431
$ws.Range("G6").Value = 821
$ws.Range("H6").Value = 130396

# Row 9: Colombia
$ws.Range("B9").Value = 702088
$ws.Range("C9").Value = 7424
$ws.Range("D9").Value = 582694
$ws.Range("E9").Value = 96876
$ws.Range("G9").Value = 243
$ws.Range("H9").Value = 22518

# Row 34: Republica Dominicana
$ws.Range("B34").Value = 102232
$ws.Range("C34").Value = 516
$ws.Range("D34").Value = 75747
$ws.Range("E34").Value = 24544
$ws.Range("G34").Value = 15
$ws.Range("H34").Value = 1941

# Row 36: Egipto
$ws.Range("B36").Value = 100708
$ws.Range("C36").Value = 151
$ws.Range("D36").Value = 82473
$ws.Range("E36").Value = 12628
$ws.Range("G36").Value = 17
$ws.Range("H36").Value = 5607

# Row 53: Barein
$ws.Range("B53").Value = 58839
$ws.Range("C53").Value = 632
$ws.Range("D53").Value = 52776
$ws.Range("E53").Value = 5855

# Row 84: Bulgaria
$ws.Range("B84").Value = 17799
$ws.Range("C84").Value = 201
$ws.Range("D84").Value = 12750
$ws.Range("E84").Value = 4336
$ws.Range("G84").Value = 7
$ws.Range("H84").Value = 713

# Row 127: Siria
$ws.Range("A127").Value = "Siria"
$ws.Range("B127").Value = 3476
$ws.Range("C127").Value = 60
$ws.Range("D127").Value = 812
$ws.Range("E127").Value = 2514
$ws.Range("G127").Value = 3
$ws.Range("H127").Value = 150

# Row 128: Tailandia
$ws.Range("A128").Value = "Tailandia"
$ws.Range("B128").Value = 3461
$ws.Range("C128").Value = 7
$ws.Range("D128").Value = 3312
$ws.Range("E128").Value = 91
$ws.Range("H128").Value = 58

# Row 135: Guadalupe
$ws.Range("A135").Value = "Guadalupe"
$ws.Range("B135").Value = 3080
$ws.Range("C135").Value = 793
$ws.Range("D135").Value = 837
$ws.Range("E135").Value = 2219
$ws.Range("H135").Value = 24

# Row 136: Jordania
$ws.Range("A136").Value = "Jordania"
$ws.Range("B136").Value = 2945
$ws.Range("C136").Value = 206
$ws.Range("D136").Value = 2084
$ws.Range("E136").Value = 840
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 21

# Row 137: Mali
$ws.Range("A137").Value = "Mali"
$ws.Range("B137").Value = 2912
$ws.Range("C137").Value = 3
$ws.Range("D137").Value = 2271
$ws.Range("E137").Value = 513
$ws.Range("H137").Value = 128

# Row 138: Aruba
$ws.Range("A138").Value = "Aruba"
$ws.Range("B138").Value = 2819
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 1407
$ws.Range("E138").Value = 1396
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 16

# Row 139: Bahamas
$ws.Range("A139").Value = "Bahamas"
$ws.Range("B139").Value = 2814
$ws.Range("C139").Value = 93
$ws.Range("D139").Value = 1220
$ws.Range("E139").Value = 1529
$ws.Range("H139").Value = 65

# Row 140: Trinidad yTobago
$ws.Range("A140").Value = "Trinidad yTobago"
$ws.Range("B140").Value = 2777
$ws.Range("C140").Value = 79
$ws.Range("D140").Value = 762
$ws.Range("E140").Value = 1970
$ws.Range("G140").Value = 2
$ws.Range("H140").Value = 45

# Row 141: Estonia
$ws.Range("A141").Value = "Estonia"
$ws.Range("B141").Value = 2632
$ws.Range("C141").Value = 32
$ws.Range("D141").Value = 2233
$ws.Range("E141").Value = 335
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 64

# Row 142: Reunion
$ws.Range("A142").Value = "Reunion"
$ws.Range("B142").Value = 2623
$ws.Range("C142").Value = 113
$ws.Range("D142").Value = 1313
$ws.Range("E142").Value = 1296
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 14

# Row 143: Sudan del Sur
$ws.Range("A143").Value = "Sudan del Sur"
$ws.Range("B143").Value = 2568
$ws.Range("C143").Value = 13
$ws.Range("D143").Value = 1290
$ws.Range("E143").Value = 1229
$ws.Range("H143").Value = 49

# Row 144: Birmania
$ws.Range("A144").Value = "Birmania"
$ws.Range("B144").Value = 2422
$ws.Range("C144").Value = 272
$ws.Range("D144").Value = 625
$ws.Range("E144").Value = 1783
$ws.Range("H144").Value = 14

# Row 163: Niger
$ws.Range("D163").Value = 1100
$ws.Range("E163").Value = 9

# Row 189: Barbados
$ws.Range("D189").Value = 158
$ws.Range("E189").Value = 15

# Row 191: Monaco
$ws.Range("B191").Value = 168
$ws.Range("C191").Value = 3
$ws.Range("D191").Value = 123
$ws.Range("E191").Value = 44

# Row 197: Islas Virgenes Britanicas
$ws.Range("B197").Value = 64
$ws.Range("C197").Value = 1
$ws.Range("D197").Value = 30
$ws.Range("E197").Value = 33
